$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp in J2
$ws.Range("J2").Value = 44401.558807870373

# GROSS SALES row (row 6)
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 65535

# Gross Returns row (row 7)
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 65535

# NET SALES row (row 8)
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 65535

# Taxes row (row 9)
$ws.Range("C9").Value = 0

# TICKET TOTAL row (row 10)
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 65535

# Cash row (row 13)
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0

# Credit Card row (row 14)
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0

# Ohio Tax row (row 21)
$ws.Range("C21").Value = 0
